$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new Price value would otherwise be parsed as a number,
# so they stay Text (matching the source data which always stores Price as text).
$ws.Range("D5,D8,D10,D11,D15,D16,D18,D20,D22,D23,D25,D26,D28,D35,D37,D38,D40,D41,D43,D45,D49").NumberFormat = "@"

$ws.Range("D2").Value = "27.601.55"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.629.63"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "212.29"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "23.37"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").Value = "0.0876"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.860.59"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "1.632.36"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "65.38"
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "27.562.35"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "231.10"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "7.57"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "10.58"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("D23").Value = "4.35"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("E24").Value = "  +8.52%  "
$ws.Range("D25").Value = "149.35"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "15.53"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "1.468.91"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D37").Value = "0.936"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("D38").Value = "0.878"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "0.556"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").Value = "1.04"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "67.88"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").Value = "1.769.40"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "87.78"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  +1.02%  "
